$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "2.0.0-sd-202406-matchbox-patch"
$meta.Range("B8").Value = "2024-06-19T17:47:42+02:00"
$meta.Range("B10").Value = "HL7 International - Structured Documents (http://www.hl7.org/Special/committees/structure, structdog@lists.HL7.org)"

# --- Elements sheet updates ---
$elem = $wb.Worksheets.Item("Elements")
$elem.Range("Z3").Value = "http://hl7.org/cda/stds/core/ValueSet/CDARoleClass"

# Update column Z width to reflect new (longer) content width, matching published output
$elem.Columns.Item(26).ColumnWidth = 49
